$wb = $excel.ActiveWorkbook
$inner = $wb.Worksheets.Item("InnerPage")

# Create HomePage right before InnerPage (created first so it receives sheetId 2,
# matching the tab order HomePage / InnerPage / LandingPage)
$home = $wb.Worksheets.Add($inner)
$home.Name = "HomePage"
$home.Range("A1").Value = "Path"
$home.Range("B1").Value = "ContentType"
$home.Range("A2").Value = "/"
$home.Range("B2").Value = "Site Home"

# Re-fetch the InnerPage reference: its Index shifted after a sheet was inserted before it
$inner2 = $wb.Worksheets.Item("InnerPage")

# Create LandingPage right after InnerPage (created second so it receives sheetId 3)
$landing = $wb.Worksheets.Add($null, $inner2)
$landing.Name = "LandingPage"
$landing.Range("A1").Value = "Path"
$landing.Range("B1").Value = "ContentType"
# Set B2 (the repeated "Landing Page" string) before A2 so new shared-string entries
# are produced in the same order as in the target workbook.
$landing.Range("B2").Value = "Landing Page"
$landing.Range("A2").Value = "/about-nci"
$landing.Range("A3").Value = "/espanol/instituto"
$landing.Range("B3").Value = "Landing Page"
$landing.Range("A4").Value = "/news-events"
$landing.Range("B4").Value = "Landing Page"
$landing.Range("A5").Value = "/espanol/noticias"
$landing.Range("B5").Value = "Landing Page"

# Restore each sheet's selection/active-cell state to match the target workbook
[void]$home.Range("A3").Select()

$inner3 = $wb.Worksheets.Item("InnerPage")
[void]$inner3.Range("B1:B1048576").Select()

[void]$landing.Range("A6").Select()
[void]$landing.Activate()
